# Add a new "campus" column header to the doubtfire_users test CSV/XLSX.
# The sheet currently has headers email..student_id in A1:G1 (row 2 holds
# one sample data row). We append a new header "campus" in H1, matching
# the commit "TEST: Update test csv files to add campus".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("doubtfire_users")

$ws.Range("H1").Value = "campus"

# Leave the new column selected, as happened in the authored edit.
[void]$ws.Range("H1").Select()
